# "edit option merge option and more"
# Adds a "##filename##" column (F) with file1/file2/file3 values, formats the
# existing C/D percentage columns (C as whole-percent, D as one-decimal percent
# via a new custom number format), widens the new column, and extends the
# sheet with an extra (still mostly blank) row that carries the D-column
# percent formatting down one more row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "##filename##" column ---------------------------------------------
$ws.Range("F1").Value = "##filename##"
$ws.Range("F2").Value = "file1"
$ws.Range("F3").Value = "file2"
$ws.Range("F4").Value = "file3"

$ws.Columns.Item(6).ColumnWidth = 11.5703125

# --- Percentage formatting on C / D -----------------------------------------
# C2:C4 -> builtin "Percent" style (whole number, e.g. 26%)
$ws.Range("C2:C4").Style = "Percent"

# D2:D5 -> "Percent" style too, but with one decimal place (e.g. 26.3%)
$ws.Range("D2:D5").Style = "Percent"
$ws.Range("D2:D5").NumberFormat = "0.0%"

# --- New trailing row ---------------------------------------------------
# Row 5 only carries D5's percent formatting (no values yet).
$ws.Range("D5").Style = "Percent"
$ws.Range("D5").NumberFormat = "0.0%"

# --- Selection / view state ----------------------------------------------
$null = $ws.Range("D6").Select()
